$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.951.68'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '3.515.46'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.08'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.94'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.67%  '
$ws.Range('D7').Value = '3.514.42'
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.480'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.80'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.06%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.423'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = '4.113.80'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '31.62'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.65%  '
$ws.Range('D16').Value = '3.522.04'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '67.123.46'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.72'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +8.71%  '
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.36'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.52%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '435.86'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.24%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.610'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '79.66'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('D25').Value = '3.654.83'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -4.07%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.85'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.37'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.42%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.58'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.89%  '
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.37'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('D35').Value = '3.510.02'
$ws.Range('E35').Value = '  +0.42%  '
$ws.Range('E36').Value = '  -3.45%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.88'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.71%  '
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '168.94'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.49%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.43'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.09'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -10.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.895'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '28.67'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '45.83'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.44'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.33%  '
$ws.Range('E51').Value = '  -0.59%  '
